# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from serial date 45500 to 45501 (i.e. advance the date by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45500) {
        $cell.Value2 = 45501
    }
}
